# Auto-generated Excel COM-interop script
# Applies numeric corrections to the Chocobo_Profits leve-crafting workbook
# across the ALC, ARM, BSM, CRP, CUL, LTW and WVR sheets (per scheduled runner update).

$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H28").Value = 596.5
$ws.Range("I28").Value = 606.15
$ws.Range("K28").Value = 606.15
$ws.Range("M28").Value = -121.15

$ws.Range("H32").Value = 800.5
$ws.Range("J32").Value = 0
$ws.Range("L32").Value = 0
$ws.Range("N32").ClearContents()

$ws.Range("H33").Value = 189.32259
$ws.Range("I33").Value = 116.5
$ws.Range("J33").Value = 367.33334
$ws.Range("K33").Value = 116.5
$ws.Range("L33").Value = 367.33334
$ws.Range("M33").Value = 112.5
$ws.Range("N33").Value = -825.33334

$ws.Range("H40").Value = 1160
$ws.Range("I40").Value = 1000
$ws.Range("K40").Value = 1000
$ws.Range("M40").Value = -825

$ws.Range("H51").Value = 5910.125
$ws.Range("I51").Value = 0
$ws.Range("K51").Value = 0
$ws.Range("M51").ClearContents()

$ws.Range("H62").Value = 4122.5
$ws.Range("I62").Value = 1993.3334
$ws.Range("K62").Value = 1993.3334
$ws.Range("M62").Value = -1369.3334

$ws.Range("H65").Value = 4122.5
$ws.Range("I65").Value = 1993.3334
$ws.Range("K65").Value = 9966.666999999999
$ws.Range("M65").Value = -6846.666999999999

$ws.Range("H86").Value = 5399.8
$ws.Range("I86").Value = 1999.5
$ws.Range("J86").Value = 7666.6665
$ws.Range("K86").Value = 1999.5
$ws.Range("L86").Value = 7666.6665
$ws.Range("M86").Value = -876.5
$ws.Range("N86").Value = -9912.666499999999

$ws.Range("H89").Value = 5399.8
$ws.Range("I89").Value = 1999.5
$ws.Range("J89").Value = 7666.6665
$ws.Range("K89").Value = 9997.5
$ws.Range("L89").Value = 38333.3325
$ws.Range("M89").Value = -4381.5
$ws.Range("N89").Value = -49565.3325

$ws.Range("H123").Value = 42183.332
$ws.Range("J123").Value = 42183.332
$ws.Range("L123").Value = 42183.332
$ws.Range("N123").Value = -51983.332

$ws.Range("H125").Value = 2618.8
$ws.Range("I125").Value = 1925
$ws.Range("J125").Value = 3081.3333
$ws.Range("K125").Value = 17325
$ws.Range("L125").Value = 27731.9997
$ws.Range("M125").Value = -14865
$ws.Range("N125").Value = -32651.9997

$ws.Range("H141").Value = 9555.933999999999
$ws.Range("I141").Value = 12153.9
$ws.Range("J141").Value = 4360
$ws.Range("K141").Value = 36461.7
$ws.Range("L141").Value = 13080
$ws.Range("M141").Value = -31281.7
$ws.Range("N141").Value = -23440

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 7418.9
$ws.Range("J32").Value = 10516.277
$ws.Range("L32").Value = 10516.277
$ws.Range("N32").Value = -11090.277

$ws.Range("H61").Value = 2093.3
$ws.Range("I61").Value = 1585.375
$ws.Range("J61").Value = 4125
$ws.Range("K61").Value = 1585.375
$ws.Range("L61").Value = 4125
$ws.Range("M61").Value = -1373.375
$ws.Range("N61").Value = -4549

$ws.Range("H74").Value = 3205.7908
$ws.Range("I74").Value = 3406.9375
$ws.Range("J74").Value = 2620.6365
$ws.Range("K74").Value = 3406.9375
$ws.Range("L74").Value = 2620.6365
$ws.Range("M74").Value = -2532.9375
$ws.Range("N74").Value = -4368.636500000001

$ws.Range("H77").Value = 3205.7908
$ws.Range("I77").Value = 3406.9375
$ws.Range("J77").Value = 2620.6365
$ws.Range("K77").Value = 17034.6875
$ws.Range("L77").Value = 13103.1825
$ws.Range("M77").Value = -12666.6875
$ws.Range("N77").Value = -21839.1825

$ws.Range("H107").Value = 0
$ws.Range("J107").Value = 0
$ws.Range("L107").Value = 0
$ws.Range("N107").ClearContents()

$ws.Range("H110").Value = 833.0833
$ws.Range("I110").Value = 750.05884
$ws.Range("J110").Value = 1034.7142
$ws.Range("K110").Value = 750.05884
$ws.Range("L110").Value = 1034.7142
$ws.Range("M110").Value = 1294.94116
$ws.Range("N110").Value = -5124.7142

$ws.Range("H136").Value = 2093.3
$ws.Range("I136").Value = 1585.375
$ws.Range("J136").Value = 4125
$ws.Range("K136").Value = 4756.125
$ws.Range("L136").Value = 12375
$ws.Range("M136").Value = -2206.125
$ws.Range("N136").Value = -17475

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H20").Value = 11517.214
$ws.Range("I20").Value = 1352.875
$ws.Range("J20").Value = 25069.666
$ws.Range("K20").Value = 1352.875
$ws.Range("L20").Value = 25069.666
$ws.Range("M20").Value = -1105.875
$ws.Range("N20").Value = -25563.666

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 11114013
$ws.Range("I31").Value = 1447.8276
$ws.Range("J31").Value = 31255538
$ws.Range("K31").Value = 1447.8276
$ws.Range("L31").Value = 31255538
$ws.Range("M31").Value = -1152.8276
$ws.Range("N31").Value = -31256128

$ws.Range("H34").Value = 11114013
$ws.Range("I34").Value = 1447.8276
$ws.Range("J34").Value = 31255538
$ws.Range("K34").Value = 1447.8276
$ws.Range("L34").Value = 31255538
$ws.Range("M34").Value = -1245.8276
$ws.Range("N34").Value = -31255942

$ws.Range("H68").Value = 99999
$ws.Range("J68").Value = 99999
$ws.Range("L68").Value = 99999
$ws.Range("N68").Value = -101497

$ws.Range("H71").Value = 99999
$ws.Range("J71").Value = 99999
$ws.Range("L71").Value = 299997
$ws.Range("N71").Value = -307485

$ws.Range("H105").Value = 2148.4167
$ws.Range("I105").Value = 2230.0588
$ws.Range("K105").Value = 2230.0588
$ws.Range("M105").Value = -483.0587999999998

$ws.Range("H138").Value = 26338.334
$ws.Range("J138").Value = 26338.334
$ws.Range("L138").Value = 26338.334
$ws.Range("N138").Value = -36618.334

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H113").Value = 576.1698
$ws.Range("I113").Value = 569.6539
$ws.Range("J113").Value = 582.44446
$ws.Range("K113").Value = 1708.9617
$ws.Range("L113").Value = 1747.33338
$ws.Range("M113").Value = 461.0382999999999
$ws.Range("N113").Value = -6087.33338

$ws.Range("H114").Value = 4687.75
$ws.Range("I114").Value = 1150
$ws.Range("J114").Value = 5395.3
$ws.Range("K114").Value = 3450
$ws.Range("L114").Value = 16185.9
$ws.Range("M114").Value = -196
$ws.Range("N114").Value = -22693.9

$ws.Range("H131").Value = 6667617.5
$ws.Range("J131").Value = 858.7143
$ws.Range("L131").Value = 2576.1429
$ws.Range("N131").Value = -12656.1429

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H55").Value = 625.53845
$ws.Range("I55").Value = 387
$ws.Range("J55").Value = 731.55554
$ws.Range("K55").Value = 387
$ws.Range("L55").Value = 731.55554
$ws.Range("M55").Value = -214
$ws.Range("N55").Value = -1077.55554

$ws.Range("H122").Value = 4857.737
$ws.Range("I122").Value = 2691
$ws.Range("J122").Value = 8572.143
$ws.Range("K122").Value = 8073
$ws.Range("L122").Value = 25716.429
$ws.Range("M122").Value = -5623
$ws.Range("N122").Value = -30616.429

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H136").Value = 1321.4
$ws.Range("I136").Value = 778
$ws.Range("J136").Value = 2241
$ws.Range("K136").Value = 2334
$ws.Range("L136").Value = 6723
$ws.Range("M136").Value = 216
$ws.Range("N136").Value = -11823

$ws.Range("H137").Value = 40770
$ws.Range("J137").Value = 48476.668
$ws.Range("L137").Value = 48476.668
$ws.Range("N137").Value = -58676.668
